$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 563.8333
$ws.Range("I15").Value = 563.8333
$ws.Range("K15").Value = 1691.4999
$ws.Range("M15").Value = -1522.4999
$ws.Range("H132").Value = 1118.15
$ws.Range("I132").Value = 1061.2106
$ws.Range("K132").Value = 3183.6318
$ws.Range("M132").Value = -653.6318000000001
$ws.Range("H137").Value = 2874.6785
$ws.Range("I137").Value = 2442.6667
$ws.Range("K137").Value = 7328.000100000001
$ws.Range("M137").Value = -4778.000100000001
$ws.Range("H138").Value = 1591082.1
$ws.Range("I138").Value = 1749.1482
$ws.Range("J138").Value = 2783082
$ws.Range("K138").Value = 5247.444600000001
$ws.Range("L138").Value = 8349246
$ws.Range("M138").Value = -107.4446000000007
$ws.Range("N138").Value = -8359526

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 359.2
$ws.Range("I5").Value = 324.25
$ws.Range("K5").Value = 324.25
$ws.Range("M5").Value = -212.25
$ws.Range("H32").Value = 4173325.2
$ws.Range("J32").Value = 33166.332
$ws.Range("L32").Value = 33166.332
$ws.Range("N32").Value = -33740.332
$ws.Range("H45").Value = 2106.7
$ws.Range("I45").Value = 2513.125
$ws.Range("J45").Value = 1835.75
$ws.Range("K45").Value = 2513.125
$ws.Range("L45").Value = 1835.75
$ws.Range("M45").Value = -2136.125
$ws.Range("N45").Value = -2589.75
$ws.Range("H61").Value = 40003196
$ws.Range("I61").Value = 1437.0667
$ws.Range("K61").Value = 1437.0667
$ws.Range("M61").Value = -1225.0667
$ws.Range("H102").Value = 2863.52
$ws.Range("I102").Value = 2338.5
$ws.Range("J102").Value = 3531.7273
$ws.Range("K102").Value = 2338.5
$ws.Range("L102").Value = 3531.7273
$ws.Range("M102").Value = -716.5
$ws.Range("N102").Value = -6775.7273
$ws.Range("H122").Value = 4598.25
$ws.Range("I122").Value = 3020.0588
$ws.Range("K122").Value = 9060.1764
$ws.Range("M122").Value = -6610.1764
$ws.Range("H135").Value = 2928713.2
$ws.Range("J135").Value = 2928713.2
$ws.Range("L135").Value = 2928713.2
$ws.Range("N135").Value = -2938853.2
$ws.Range("H136").Value = 40003196
$ws.Range("I136").Value = 1437.0667
$ws.Range("K136").Value = 4311.2001
$ws.Range("M136").Value = -1761.2001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 359.2
$ws.Range("I4").Value = 324.25
$ws.Range("K4").Value = 324.25
$ws.Range("M4").Value = -209.25
$ws.Range("H20").Value = 5379994.5
$ws.Range("I20").Value = 7579094.5
$ws.Range("K20").Value = 7579094.5
$ws.Range("M20").Value = -7578847.5
$ws.Range("H99").Value = 5052419
$ws.Range("I99").Value = 1268.6364
$ws.Range("K99").Value = 1268.6364
$ws.Range("M99").Value = 229.3635999999999
$ws.Range("H105").Value = 3869.375
$ws.Range("I105").Value = 2962.375
$ws.Range("K105").Value = 2962.375
$ws.Range("M105").Value = -1215.375
$ws.Range("H107").Value = 62515850
$ws.Range("I107").Value = 112517370
$ws.Range("J107").Value = 13947.75
$ws.Range("K107").Value = 112517370
$ws.Range("L107").Value = 13947.75
$ws.Range("M107").Value = -112515450
$ws.Range("N107").Value = -17787.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4646.6924
$ws.Range("J31").Value = 6986.2856
$ws.Range("L31").Value = 6986.2856
$ws.Range("N31").Value = -7576.2856
$ws.Range("H34").Value = 4646.6924
$ws.Range("J34").Value = 6986.2856
$ws.Range("L34").Value = 6986.2856
$ws.Range("N34").Value = -7390.2856
$ws.Range("H58").Value = 5164.615
$ws.Range("I58").Value = 2128.625
$ws.Range("J58").Value = 6513.9443
$ws.Range("K58").Value = 2128.625
$ws.Range("L58").Value = 6513.9443
$ws.Range("M58").Value = -1925.625
$ws.Range("N58").Value = -6919.9443
$ws.Range("H99").Value = 3412.7
$ws.Range("I99").Value = 1737.091
$ws.Range("J99").Value = 5460.6665
$ws.Range("K99").Value = 1737.091
$ws.Range("L99").Value = 5460.6665
$ws.Range("M99").Value = -239.0909999999999
$ws.Range("N99").Value = -8456.666499999999
$ws.Range("H126").Value = 3412.7
$ws.Range("I126").Value = 1737.091
$ws.Range("J126").Value = 5460.6665
$ws.Range("K126").Value = 5211.272999999999
$ws.Range("L126").Value = 16381.9995
$ws.Range("M126").Value = -2741.272999999999
$ws.Range("N126").Value = -21321.9995
$ws.Range("H136").Value = 5164.615
$ws.Range("I136").Value = 2128.625
$ws.Range("J136").Value = 6513.9443
$ws.Range("K136").Value = 6385.875
$ws.Range("L136").Value = 19541.8329
$ws.Range("M136").Value = -3835.875
$ws.Range("N136").Value = -24641.8329

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 55555996
$ws.Range("I98").Value = 325
$ws.Range("J98").Value = 83333830
$ws.Range("K98").Value = 975
$ws.Range("L98").Value = 250001490
$ws.Range("M98").Value = 523
$ws.Range("N98").Value = -250004486
$ws.Range("H141").Value = 3043.2307
$ws.Range("I141").Value = 3043.2307
$ws.Range("K141").Value = 9129.6921
$ws.Range("M141").Value = -3949.6921

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7434
$ws.Range("I70").Value = 4859.778
$ws.Range("K70").Value = 4859.778
$ws.Range("M70").Value = -4589.778
$ws.Range("H73").Value = 7434
$ws.Range("I73").Value = 4859.778
$ws.Range("K73").Value = 4859.778
$ws.Range("M73").Value = -3923.778
$ws.Range("H80").Value = 4000.6
$ws.Range("I80").Value = 1686
$ws.Range("J80").Value = 7472.5
$ws.Range("K80").Value = 1686
$ws.Range("L80").Value = 7472.5
$ws.Range("M80").Value = -688
$ws.Range("N80").Value = -9468.5
$ws.Range("H83").Value = 4000.6
$ws.Range("I83").Value = 1686
$ws.Range("J83").Value = 7472.5
$ws.Range("K83").Value = 8430
$ws.Range("L83").Value = 37362.5
$ws.Range("M83").Value = -3438
$ws.Range("N83").Value = -47346.5
$ws.Range("H113").Value = 5802.109
$ws.Range("I113").Value = 2176.087
$ws.Range("K113").Value = 2176.087
$ws.Range("M113").Value = -6.086999999999989

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1925000
$ws.Range("I20").Value = 2550000
$ws.Range("K20").Value = 2550000
$ws.Range("M20").Value = -2549774
$ws.Range("H22").Value = 2085.4285
$ws.Range("I22").Value = 498.5
$ws.Range("J22").Value = 2720.2
$ws.Range("K22").Value = 498.5
$ws.Range("L22").Value = 2720.2
$ws.Range("M22").Value = -203.5
$ws.Range("N22").Value = -3310.2
$ws.Range("H27").Value = 2085.4285
$ws.Range("I27").Value = 498.5
$ws.Range("J27").Value = 2720.2
$ws.Range("K27").Value = 498.5
$ws.Range("L27").Value = 2720.2
$ws.Range("M27").Value = -391.5
$ws.Range("N27").Value = -2934.2
$ws.Range("H40").Value = 33337552
$ws.Range("I40").Value = 50002330
$ws.Range("K40").Value = 50002330
$ws.Range("M40").Value = -50002194
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -52246
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -161232
$ws.Range("H100").Value = 3409.3333
$ws.Range("I100").Value = 2414.2856
$ws.Range("J100").Value = 5399.4287
$ws.Range("K100").Value = 2414.2856
$ws.Range("L100").Value = 5399.4287
$ws.Range("M100").Value = -1873.2856
$ws.Range("N100").Value = -6481.4287
$ws.Range("H122").Value = 4261.9585
$ws.Range("I122").Value = 3378.4167
$ws.Range("K122").Value = 10135.2501
$ws.Range("M122").Value = -7685.250100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6938.6665
$ws.Range("J62").Value = 6824.143
$ws.Range("L62").Value = 6824.143
$ws.Range("N62").Value = -8072.143
$ws.Range("H65").Value = 6938.6665
$ws.Range("J65").Value = 6824.143
$ws.Range("L65").Value = 34120.715
$ws.Range("N65").Value = -40360.715
